$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old "Hello World" placeholder rows; the new layout starts on row 3.
$ws.Range("A1").ClearContents()
$ws.Range("A2").ClearContents()

$ws.Range("A3").Value = "FİYAT TEKLİFİ"
$ws.Range("A4").Value = "FİRMA ADI"
$ws.Range("A5").Value = "YETKİLİ ADI"
$ws.Range("A6").Value = "TELEFON"
$ws.Range("A7").Value = "E-POSTA"
$ws.Range("A8").Value = "FATURA ADRESİ"
$ws.Range("A9").Value = "VERGİ DAİRESİ/NO"

$ws.Range("A10").Value = "S.NO"
$ws.Range("B10").Value = "ÜRÜN ADI "
$ws.Range("C10").Value = "MODEL "
$ws.Range("D10").Value = "ÖLÇÜ "
$ws.Range("E10").Value = "RENK"
$ws.Range("F10").Value = "MİKTAR"
$ws.Range("G10").Value = "BİRİM FİYATI "
$ws.Range("H10").Value = "TUTAR"
$ws.Range("I10").Value = "GÖRSEL "
